$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.442.90'
$ws.Range("E2").Value = '  +0.42%  '
$ws.Range("D3").Value = '1.635.65'
$ws.Range("E3").Value = '  -0.77%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '''212.23'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.64%  '
$ws.Range("D6").Value = '''0.531'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.53%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '''22.85'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.53%  '
$ws.Range("E9").Value = '  -2.17%  '
$ws.Range("E10").Value = '  -0.76%  '
$ws.Range("E11").Value = '  +1.27%  '
$ws.Range("D12").Value = '1.868.42'
$ws.Range("E12").Value = '  -0.66%  '
$ws.Range("D13").Value = '1.632.70'
$ws.Range("E13").Value = '  -0.97%  '
$ws.Range("D14").Value = '''0.569'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.43%  '
$ws.Range("E15").Value = '  -1.63%  '
$ws.Range("D16").Value = '''64.11'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.41%  '
$ws.Range("D17").Value = '27.473.20'
$ws.Range("E17").Value = '  +0.52%  '
$ws.Range("D18").Value = '''227.73'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.67%  '
$ws.Range("E19").Value = '  -0.05%  '
$ws.Range("D20").Value = '''7.64'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.98%  '
$ws.Range("E21").Value = '  -0.06%  '
$ws.Range("E22").Value = '  -2.07%  '
$ws.Range("D23").Value = '''9.77'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +5.59%  '
$ws.Range("E24").Value = '  -2.98%  '
$ws.Range("D25").Value = '''149.87'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.62%  '
$ws.Range("D26").Value = '''6.96'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.78%  '
$ws.Range("E27").Value = '  +1.68%  '
$ws.Range("E28").Value = '  -0.08%  '
$ws.Range("D29").Value = '''15.54'
$ws.Range("D29").Style = "Normal"
$ws.Range("E30").Value = '  -0.72%  '
$ws.Range("E31").Value = '  -1.82%  '
$ws.Range("D32").Value = '''3.28'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.44%  '
$ws.Range("E33").Value = '  +3.10%  '
$ws.Range("D34").Value = '1.408.60'
$ws.Range("E34").Value = '  -3.34%  '
$ws.Range("D35").Value = '''1.58'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.46%  '
$ws.Range("E36").Value = '  -2.23%  '
$ws.Range("D37").Value = '''0.569'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.22%  '
$ws.Range("D38").Value = '''0.0167'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.12%  '
$ws.Range("B39").Value = 'ARBITRUM'
$ws.Range("C39").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D39").Value = '''0.870'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.96%  '
$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").Value = '''0.915'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +17.43%  '
$ws.Range("E41").Value = '  -0.69%  '
$ws.Range("E42").Value = '  -0.03%  '
$ws.Range("D43").Value = '''5.50'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.63%  '
$ws.Range("D44").Value = '''2.24'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.36%  '
$ws.Range("D45").Value = '''64.64'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.80%  '
$ws.Range("D46").Value = '1.777.48'
$ws.Range("E46").Value = '  -0.61%  '
$ws.Range("E47").Value = '  -3.20%  '
$ws.Range("D48").Value = '''85.83'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.59%  '
$ws.Range("D49").Value = '0.0₆0106'
$ws.Range("E49").Value = '  -0.36%  '
$ws.Range("D50").Value = '''0.0986'
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Value = '''7.72'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.65%  '
